$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.636.46'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '1.737.94'
$ws.Range("E3").Value = '  -1.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +1.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.92'
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3828'
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3342'
$ws.Range("E8").Value = '  -2.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.44'
$ws.Range("E9").Value = '  -5.47%  '

$ws.Range("E10").Value = '  -3.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07145'
$ws.Range("E11").Value = '  -3.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.00'
$ws.Range("E13").Value = '  -1.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.091'
$ws.Range("E14").Value = '  -4.06%  '

$ws.Range("D15").Value = '1.752.62'
$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.960'
$ws.Range("E16").Value = '  -1.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001044'
$ws.Range("E17").Value = '  -2.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06579'
$ws.Range("E18").Value = '  -1.56%  '

$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.43'
$ws.Range("E20").Value = '  -4.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.60'
$ws.Range("E21").Value = '  -4.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.147'
$ws.Range("E22").Value = '  -4.16%  '

$ws.Range("D23").Value = '27.704.73'
$ws.Range("E23").Value = '  -0.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.45'
$ws.Range("E24").Value = '  -4.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.395'
$ws.Range("E25").Value = '  +0.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.66'
$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.60'
$ws.Range("E27").Value = '  -5.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.250'
$ws.Range("E28").Value = '  -6.87%  '

$ws.Range("D29").Value = '1.943.15'
$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.264'
$ws.Range("E30").Value = '  -12.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.18'
$ws.Range("E31").Value = '  -3.62%  '

$ws.Range("E32").Value = '  +1.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.727'
$ws.Range("E33").Value = '  -6.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08589'
$ws.Range("E34").Value = '  -2.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.88'
$ws.Range("E35").Value = '  -6.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.529'
$ws.Range("E36").Value = '  +0.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.093'
$ws.Range("E37").Value = '  -4.18%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02251'
$ws.Range("E38").Value = '  -7.26%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6420'
$ws.Range("E39").Value = '  -6.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06020'
$ws.Range("E40").Value = '  -4.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2077'
$ws.Range("E41").Value = '  -4.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.188'
$ws.Range("E42").Value = '  -3.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.900'
$ws.Range("E44").Value = '  -4.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.57'
$ws.Range("E45").Value = '  -3.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.803'
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5938'
$ws.Range("E47").Value = '  -5.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.50'
$ws.Range("E48").Value = '  -4.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.968'
$ws.Range("E49").Value = '  -5.33%  '

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.142'
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06920'
$ws.Range("E51").Value = '  -6.17%  '
